$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.890.90'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.544.79'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''206.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '''0.486'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = '''0.246'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''21.27'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").Value = '''0.0582'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = '1.766.07'
$ws.Range("D13").Value = '1.547.88'
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("D14").Value = '''3.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '26.881.49'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''61.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '''213.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '0.0₃0682'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '''7.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").Value = '''9.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").Value = '''151.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''6.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = '''14.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = '''0.0458'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").Value = '''3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("D33").Value = '1.357.19'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("D34").Value = '''2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").Value = '''1.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '''0.957'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").Value = '''0.519'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '''0.802'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = '''5.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.21%  '
$ws.Range("D43").Value = '''0.990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("E44").Value = '  +2.01%  '
$ws.Range("D45").Value = '''63.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = '1.680.63'
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("D48").Value = '''85.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("D51").Value = '''0.0946'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.02%  '
